# nowe faktury z Manexu
# Add three new invoice rows to the "Faktury Manex" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Faktury Manex")

# Row 5 - invoice 4877/T/08/2013, "Rura woda"
$ws.Range("A5").Value = 41498
$ws.Range("B5").Value = 41501
$ws.Range("C5").Value = "4877/T/08/2013"
$ws.Range("D5").Value = 19.08
$ws.Range("E5").Value = "Rura woda"

# Row 6 - invoice 4916/T/08/2013, "Folia"
$ws.Range("A6").Value = 41498
$ws.Range("B6").Value = 41501
$ws.Range("C6").Value = "4916/T/08/2013"
$ws.Range("D6").Value = 115.01
$ws.Range("E6").Value = "Folia"

# Row 7 - invoice 4969/T/08/2013, "Styropian"
$ws.Range("A7").Value = 41499
$ws.Range("B7").Value = 41502
$ws.Range("C7").Value = "4969/T/08/2013"
$ws.Range("D7").Value = 8280.36
$ws.Range("E7").Value = "Styropian"

# Match the red "unpaid" formatting used on row 4 (dates s=12, rest s=13)
$ws.Range("A4:E4").Copy() | Out-Null
$ws.Range("A5:E7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("F6").Select() | Out-Null
